# Apply "Natmi following Dr Hou advice" edit: expand Bmp7-Acvr2b LR-pair table
# from a single sending cluster (sCs) to two sending clusters (Neutro, sCs),
# each paired against all 6 target clusters (ECs, FAPs, M1, M2, Neutro, sCs).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# then the 16 numeric metric columns E..T
$rows = @(
    @("Neutro", "Bmp7", "Acvr2b", "ECs", 1, 0.3333333333333333, 0.1301303333333333, 0.390391, 0.4028730131292778, 0.5029891437391611, 2, 1, 1.921801, 3.843602, 0.2963026861879742, 0.2443427896918877, 0.2500846047303333, 1.500507628382, 0.119372355982848, 0.1229017705659605),
    @("Neutro", "Bmp7", "Acvr2b", "FAPs", 1, 0.3333333333333333, 0.1301303333333333, 0.390391, 0.4028730131292778, 0.5029891437391611, 3, 1, 0.1523196666666667, 0.456959, 0.02348459928608122, 0.02904947932559493, 0.01982140899655555, 0.178392680969, 0.009461311276517226, 0.01461157273204946),
    @("Neutro", "Bmp7", "Acvr2b", "M1", 1, 0.3333333333333333, 0.1301303333333333, 0.390391, 0.4028730131292778, 0.5029891437391611, 2, 0.6666666666666666, 0.4109286666666667, 1.232786, 0.06335685524410488, 0.07836981309019599, 0.05347428436955555, 0.4812685593259999, 0.02552476717458802, 0.03941916518123578),
    @("Neutro", "Bmp7", "Acvr2b", "M2", 1, 0.3333333333333333, 0.1301303333333333, 0.390391, 0.4028730131292778, 0.5029891437391611, 3, 1, 0.9978816666666667, 2.993645, 0.1538530879789666, 0.1903099151907953, 0.1298546739105556, 1.168692065195, 0.06198325713333015, 0.0957238212868905),
    @("Neutro", "Bmp7", "Acvr2b", "Neutro", 1, 0.3333333333333333, 0.1301303333333333, 0.390391, 0.4028730131292778, 0.5029891437391611, 3, 1, 1.197361, 3.592083, 0.1846087501446398, 0.2283533989796043, 0.1558129860503333, 1.402316874453, 0.07437388342080105, 0.1148592806226782),
    @("Neutro", "Bmp7", "Acvr2b", "sCs", 1, 0.3333333333333333, 0.1301303333333333, 0.390391, 0.4028730131292778, 0.5029891437391611, 2, 1, 1.8056465, 3.611293, 0.2783940211582333, 0.2295746037219218, 0.2349693809271666, 1.409816285563, 0.1121574381411934, 0.1154735333503467),
    @("sCs", "Bmp7", "Acvr2b", "ECs", 2, 1, 0.1928755, 0.385751, 0.5971269868707222, 0.4970108562608389, 2, 1, 1.921801, 3.843602, 0.2963026861879742, 0.2443427896918877, 0.3706683287755, 1.482673315102, 0.1769303302051262, 0.1214410191259272),
    @("sCs", "Bmp7", "Acvr2b", "FAPs", 2, 1, 0.1928755, 0.385751, 0.5971269868707222, 0.4970108562608389, 3, 1, 0.1523196666666667, 0.456959, 0.02348459928608122, 0.02904947932559493, 0.02937873186816667, 0.176272391209, 0.01402328800956399, 0.01443790659354547),
    @("sCs", "Bmp7", "Acvr2b", "M1", 2, 1, 0.1928755, 0.385751, 0.5971269868707222, 0.4970108562608389, 2, 0.6666666666666666, 0.4109286666666667, 1.232786, 0.06335685524410488, 0.07836981309019599, 0.07925807204766668, 0.475548432286, 0.03783208806951686, 0.03895064790896021),
    @("sCs", "Bmp7", "Acvr2b", "M2", 2, 1, 0.1928755, 0.385751, 0.5971269868707222, 0.4970108562608389, 3, 1, 0.9978816666666667, 2.993645, 0.1538530879789666, 0.1903099151907953, 0.1924669253991667, 1.154801552395, 0.09186983084563645, 0.09458609390390481),
    @("sCs", "Bmp7", "Acvr2b", "Neutro", 2, 1, 0.1928755, 0.385751, 0.5971269868707222, 0.4970108562608389, 3, 1, 1.197361, 3.592083, 0.1846087501446398, 0.2283533989796043, 0.2309416015555, 1.385649609333, 0.1102348667238388, 0.1134941183569261),
    @("sCs", "Bmp7", "Acvr2b", "sCs", 2, 1, 0.1928755, 0.385751, 0.5971269868707222, 0.4970108562608389, 2, 1, 1.8056465, 3.611293, 0.2783940211582333, 0.2295746037219218, 0.34826497151075, 1.393059886043, 0.1662365830170399, 0.1141010703715751),
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $rowVals = $rows[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $rowVals[$j]
    }
}

Write-Output "Updated $($rows.Count) data rows (rows 2-13)"
